$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update/extend the data table (rows 2-13) with new TPM-derived values.
# Row 2: ECs -> FAPs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Nlgn3"
$ws.Cells.Item(2, 3).Value = "Nrxn1"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.1751946666666667
$ws.Cells.Item(2, 8).Value = 0.525584
$ws.Cells.Item(2, 9).Value = 0.07933130249481599
$ws.Cells.Item(2, 10).Value = 0.079331302494816
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.2191816666666667
$ws.Cells.Item(2, 14).Value = 0.657545
$ws.Cells.Item(2, 15).Value = 0.8459226744507667
$ws.Cells.Item(2, 16).Value = 0.8459226744507669
$ws.Cells.Item(2, 17).Value = 0.03839945903111112
$ws.Cells.Item(2, 18).Value = 0.3455951312800001
$ws.Cells.Item(2, 19).Value = 0.06710814757407753
$ws.Cells.Item(2, 20).Value = 0.06710814757407754

# Row 3: ECs -> MuSCs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Nlgn3"
$ws.Cells.Item(3, 3).Value = "Nrxn1"
$ws.Cells.Item(3, 4).Value = "MuSCs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.1751946666666667
$ws.Cells.Item(3, 8).Value = 0.525584
$ws.Cells.Item(3, 9).Value = 0.07933130249481599
$ws.Cells.Item(3, 10).Value = 0.079331302494816
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.03648100000000001
$ws.Cells.Item(3, 14).Value = 0.109443
$ws.Cells.Item(3, 15).Value = 0.1407969268413801
$ws.Cells.Item(3, 16).Value = 0.1407969268413801
$ws.Cells.Item(3, 17).Value = 0.006391276634666669
$ws.Cells.Item(3, 18).Value = 0.05752148971200001
$ws.Cells.Item(3, 19).Value = 0.011169603593594
$ws.Cells.Item(3, 20).Value = 0.011169603593594

# Row 4: ECs -> Resolving-Mac
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Nlgn3"
$ws.Cells.Item(4, 3).Value = "Nrxn1"
$ws.Cells.Item(4, 4).Value = "Resolving-Mac"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.1751946666666667
$ws.Cells.Item(4, 8).Value = 0.525584
$ws.Cells.Item(4, 9).Value = 0.07933130249481599
$ws.Cells.Item(4, 10).Value = 0.079331302494816
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.003441
$ws.Cells.Item(4, 14).Value = 0.010323
$ws.Cells.Item(4, 15).Value = 0.0132803987078531
$ws.Cells.Item(4, 16).Value = 0.0132803987078531
$ws.Cells.Item(4, 17).Value = 0.0006028448480000001
$ws.Cells.Item(4, 18).Value = 0.005425603632000001
$ws.Cells.Item(4, 19).Value = 0.001053551327144457
$ws.Cells.Item(4, 20).Value = 0.001053551327144458

# Row 5: FAPs -> FAPs
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Nlgn3"
$ws.Cells.Item(5, 3).Value = "Nrxn1"
$ws.Cells.Item(5, 4).Value = "FAPs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.8048609999999999
$ws.Cells.Item(5, 8).Value = 2.414583
$ws.Cells.Item(5, 9).Value = 0.364455566325916
$ws.Cells.Item(5, 10).Value = 0.3644555663259161
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.2191816666666667
$ws.Cells.Item(5, 14).Value = 0.657545
$ws.Cells.Item(5, 15).Value = 0.8459226744507667
$ws.Cells.Item(5, 16).Value = 0.8459226744507669
$ws.Cells.Item(5, 17).Value = 0.176410775415
$ws.Cells.Item(5, 18).Value = 1.587696978735
$ws.Cells.Item(5, 19).Value = 0.3083012273848876
$ws.Cells.Item(5, 20).Value = 0.3083012273848878

# Row 6: FAPs -> MuSCs
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Nlgn3"
$ws.Cells.Item(6, 3).Value = "Nrxn1"
$ws.Cells.Item(6, 4).Value = "MuSCs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.8048609999999999
$ws.Cells.Item(6, 8).Value = 2.414583
$ws.Cells.Item(6, 9).Value = 0.364455566325916
$ws.Cells.Item(6, 10).Value = 0.3644555663259161
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.03648100000000001
$ws.Cells.Item(6, 14).Value = 0.109443
$ws.Cells.Item(6, 15).Value = 0.1407969268413801
$ws.Cells.Item(6, 16).Value = 0.1407969268413801
$ws.Cells.Item(6, 17).Value = 0.029362134141
$ws.Cells.Item(6, 18).Value = 0.264259207269
$ws.Cells.Item(6, 19).Value = 0.05131422370892374
$ws.Cells.Item(6, 20).Value = 0.05131422370892377

# Row 7: FAPs -> Resolving-Mac
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Nlgn3"
$ws.Cells.Item(7, 3).Value = "Nrxn1"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.8048609999999999
$ws.Cells.Item(7, 8).Value = 2.414583
$ws.Cells.Item(7, 9).Value = 0.364455566325916
$ws.Cells.Item(7, 10).Value = 0.3644555663259161
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.003441
$ws.Cells.Item(7, 14).Value = 0.010323
$ws.Cells.Item(7, 15).Value = 0.0132803987078531
$ws.Cells.Item(7, 16).Value = 0.0132803987078531
$ws.Cells.Item(7, 17).Value = 0.002769526701
$ws.Cells.Item(7, 18).Value = 0.024925740309
$ws.Cells.Item(7, 19).Value = 0.004840115232104563
$ws.Cells.Item(7, 20).Value = 0.004840115232104565

# Row 8: MuSCs -> FAPs
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Nlgn3"
$ws.Cells.Item(8, 3).Value = "Nrxn1"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.142989
$ws.Cells.Item(8, 8).Value = 3.428967
$ws.Cells.Item(8, 9).Value = 0.5175660186035755
$ws.Cells.Item(8, 10).Value = 0.5175660186035757
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.2191816666666667
$ws.Cells.Item(8, 14).Value = 0.657545
$ws.Cells.Item(8, 15).Value = 0.8459226744507667
$ws.Cells.Item(8, 16).Value = 0.8459226744507669
$ws.Cells.Item(8, 17).Value = 0.2505222340016667
$ws.Cells.Item(8, 18).Value = 2.254700106015
$ws.Cells.Item(8, 19).Value = 0.4378208306619719
$ws.Cells.Item(8, 20).Value = 0.4378208306619721

# Row 9: MuSCs -> MuSCs
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Nlgn3"
$ws.Cells.Item(9, 3).Value = "Nrxn1"
$ws.Cells.Item(9, 4).Value = "MuSCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.142989
$ws.Cells.Item(9, 8).Value = 3.428967
$ws.Cells.Item(9, 9).Value = 0.5175660186035755
$ws.Cells.Item(9, 10).Value = 0.5175660186035757
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.03648100000000001
$ws.Cells.Item(9, 14).Value = 0.109443
$ws.Cells.Item(9, 15).Value = 0.1407969268413801
$ws.Cells.Item(9, 16).Value = 0.1407969268413801
$ws.Cells.Item(9, 17).Value = 0.04169738170900001
$ws.Cells.Item(9, 18).Value = 0.375276435381
$ws.Cells.Item(9, 19).Value = 0.072871704856912
$ws.Cells.Item(9, 20).Value = 0.07287170485691202

# Row 10: MuSCs -> Resolving-Mac
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Nlgn3"
$ws.Cells.Item(10, 3).Value = "Nrxn1"
$ws.Cells.Item(10, 4).Value = "Resolving-Mac"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.142989
$ws.Cells.Item(10, 8).Value = 3.428967
$ws.Cells.Item(10, 9).Value = 0.5175660186035755
$ws.Cells.Item(10, 10).Value = 0.5175660186035757
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.003441
$ws.Cells.Item(10, 14).Value = 0.010323
$ws.Cells.Item(10, 15).Value = 0.0132803987078531
$ws.Cells.Item(10, 16).Value = 0.0132803987078531
$ws.Cells.Item(10, 17).Value = 0.003933025149
$ws.Cells.Item(10, 18).Value = 0.035397226341
$ws.Cells.Item(10, 19).Value = 0.006873483084691596
$ws.Cells.Item(10, 20).Value = 0.0068734830846916

# Row 11: Resolving-Mac -> FAPs
$ws.Cells.Item(11, 1).Value = "Resolving-Mac"
$ws.Cells.Item(11, 2).Value = "Nlgn3"
$ws.Cells.Item(11, 3).Value = "Nrxn1"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.085348
$ws.Cells.Item(11, 8).Value = 0.256044
$ws.Cells.Item(11, 9).Value = 0.0386471125756923
$ws.Cells.Item(11, 10).Value = 0.0386471125756923
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.2191816666666667
$ws.Cells.Item(11, 14).Value = 0.657545
$ws.Cells.Item(11, 15).Value = 0.8459226744507667
$ws.Cells.Item(11, 16).Value = 0.8459226744507669
$ws.Cells.Item(11, 17).Value = 0.01870671688666667
$ws.Cells.Item(11, 18).Value = 0.16836045198
$ws.Cells.Item(11, 19).Value = 0.03269246882982949
$ws.Cells.Item(11, 20).Value = 0.0326924688298295

# Row 12: Resolving-Mac -> MuSCs
$ws.Cells.Item(12, 1).Value = "Resolving-Mac"
$ws.Cells.Item(12, 2).Value = "Nlgn3"
$ws.Cells.Item(12, 3).Value = "Nrxn1"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.085348
$ws.Cells.Item(12, 8).Value = 0.256044
$ws.Cells.Item(12, 9).Value = 0.0386471125756923
$ws.Cells.Item(12, 10).Value = 0.0386471125756923
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 0.03648100000000001
$ws.Cells.Item(12, 14).Value = 0.109443
$ws.Cells.Item(12, 15).Value = 0.1407969268413801
$ws.Cells.Item(12, 16).Value = 0.1407969268413801
$ws.Cells.Item(12, 17).Value = 0.003113580388000001
$ws.Cells.Item(12, 18).Value = 0.028022223492
$ws.Cells.Item(12, 19).Value = 0.005441394681950329
$ws.Cells.Item(12, 20).Value = 0.005441394681950331

# Row 13: Resolving-Mac -> Resolving-Mac
$ws.Cells.Item(13, 1).Value = "Resolving-Mac"
$ws.Cells.Item(13, 2).Value = "Nlgn3"
$ws.Cells.Item(13, 3).Value = "Nrxn1"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.085348
$ws.Cells.Item(13, 8).Value = 0.256044
$ws.Cells.Item(13, 9).Value = 0.0386471125756923
$ws.Cells.Item(13, 10).Value = 0.0386471125756923
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.003441
$ws.Cells.Item(13, 14).Value = 0.010323
$ws.Cells.Item(13, 15).Value = 0.0132803987078531
$ws.Cells.Item(13, 16).Value = 0.0132803987078531
$ws.Cells.Item(13, 17).Value = 0.000293682468
$ws.Cells.Item(13, 18).Value = 0.002643142212
$ws.Cells.Item(13, 19).Value = 0.0005132490639124771
$ws.Cells.Item(13, 20).Value = 0.0005132490639124773
